# Update the generated "Installed Asset IDs" tokens (Vendor_<timestamp>)
# in column A of each per-vendor worksheet, matching a re-run of the
# installed-software report generator.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "Microsoft";              Cell = "A2"; Value = "Microsoft_1663605558641" },
    @{ Sheet = "CISCO";                  Cell = "A2"; Value = "CISCO_1663605558782" },
    @{ Sheet = "Fortinet Technologies";  Cell = "A2"; Value = "Fortinet_1663605558910" },
    @{ Sheet = "Citrix";                 Cell = "A2"; Value = "Citrix_1663605558995" },
    @{ Sheet = "Citrix";                 Cell = "A3"; Value = "Citrix_1663605559148" },
    @{ Sheet = "Citrix";                 Cell = "A4"; Value = "Citrix_1663605559236" },
    @{ Sheet = "Microsoft";              Cell = "A3"; Value = "Microsoft_1663605559303" },
    @{ Sheet = "Oracle";                 Cell = "A2"; Value = "Oracle_1663605559389" },
    @{ Sheet = "Vendor Not Found";       Cell = "A2"; Value = "Vendor_1663605559563" },
    @{ Sheet = "Test Vendor 1";          Cell = "A2"; Value = "Test_1663605559611" },
    @{ Sheet = "Test Vendor 2";          Cell = "A2"; Value = "Test_1663605559732" },
    @{ Sheet = "Waves Audio";            Cell = "A2"; Value = "Waves_1663605559984" },
    @{ Sheet = "Citrix";                 Cell = "A5"; Value = "Citrix_1663605560032" },
    @{ Sheet = "Microsoft";              Cell = "A4"; Value = "Microsoft_1663605560077" },
    @{ Sheet = "Oracle";                 Cell = "A3"; Value = "Oracle_1663605560197" },
    @{ Sheet = "CheckPoint";             Cell = "A2"; Value = "CheckPoint_1663605560240" },
    @{ Sheet = "Adobe";                  Cell = "A2"; Value = "Adobe_1663605560295" },
    @{ Sheet = "Citrix";                 Cell = "A6"; Value = "Citrix_1663605560333" },
    @{ Sheet = "CheckPoint";             Cell = "A3"; Value = "CheckPoint_1663605560437" },
    @{ Sheet = "Business Objects";       Cell = "A2"; Value = "Business_1663605560472" },
    @{ Sheet = "ConnectWise";            Cell = "A2"; Value = "ConnectWise_1663605560672" },
    @{ Sheet = "CISCO";                  Cell = "A3"; Value = "CISCO_1663605560710" },
    @{ Sheet = "SAP";                    Cell = "A2"; Value = "SAP_1663605560866" }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
}
